$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 499
$ws.Range("J32").Value = 499
$ws.Range("L32").Value = 499
$ws.Range("N32").Value = -1151

$ws.Range("H38").Value = 2177.8572
$ws.Range("I38").Value = 2177.8572
$ws.Range("J38").Value = 0
$ws.Range("K38").Value = 6533.571599999999
$ws.Range("L38").Value = 0
$ws.Range("M38").Value = -6161.571599999999
$ws.Range("N38").Value = $null

$ws.Range("H80").Value = 2074.5334
$ws.Range("I80").Value = 1640.2858
$ws.Range("J80").Value = 2454.5
$ws.Range("K80").Value = 4920.857400000001
$ws.Range("L80").Value = 7363.5
$ws.Range("M80").Value = -3922.857400000001
$ws.Range("N80").Value = -9359.5

$ws.Range("H83").Value = 2074.5334
$ws.Range("I83").Value = 1640.2858
$ws.Range("J83").Value = 2454.5
$ws.Range("K83").Value = 14762.5722
$ws.Range("L83").Value = 22090.5
$ws.Range("M83").Value = -9770.572200000001
$ws.Range("N83").Value = -32074.5

$ws.Range("H95").Value = 48331.668
$ws.Range("J95").Value = 48331.668
$ws.Range("L95").Value = 48331.668
$ws.Range("N95").Value = -53823.668

$ws.Range("H116").Value = 7667
$ws.Range("I116").Value = 8000
$ws.Range("J116").Value = 7001
$ws.Range("K116").Value = 8000
$ws.Range("L116").Value = 7001
$ws.Range("M116").Value = -4558
$ws.Range("N116").Value = -13885

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 416.33334
$ws.Range("I2").Value = 416.33334
$ws.Range("K2").Value = 416.33334
$ws.Range("M2").Value = -303.33334

$ws.Range("H116").Value = 416.33334
$ws.Range("I116").Value = 416.33334
$ws.Range("K116").Value = 416.33334
$ws.Range("M116").Value = 1877.66666

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 416.33334
$ws.Range("I3").Value = 416.33334
$ws.Range("K3").Value = 416.33334
$ws.Range("M3").Value = -302.33334

$ws.Range("H20").Value = 1966.7273
$ws.Range("I20").Value = 1427.1111
$ws.Range("J20").Value = 4395
$ws.Range("K20").Value = 1427.1111
$ws.Range("L20").Value = 4395
$ws.Range("M20").Value = -1180.1111
$ws.Range("N20").Value = -4889

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H25").Value = 475.0909
$ws.Range("I25").Value = 500
$ws.Range("J25").Value = 408.66666
$ws.Range("K25").Value = 500
$ws.Range("L25").Value = 408.66666
$ws.Range("M25").Value = -326
$ws.Range("N25").Value = -756.66666

$ws.Range("H58").Value = 5248.143
$ws.Range("I58").Value = 1996.3334
$ws.Range("J58").Value = 7687
$ws.Range("K58").Value = 1996.3334
$ws.Range("L58").Value = 7687
$ws.Range("M58").Value = -1793.3334
$ws.Range("N58").Value = -8093

$ws.Range("H136").Value = 5248.143
$ws.Range("I136").Value = 1996.3334
$ws.Range("J136").Value = 7687
$ws.Range("K136").Value = 5989.0002
$ws.Range("L136").Value = 23061
$ws.Range("M136").Value = -3439.0002
$ws.Range("N136").Value = -28161

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H41").Value = 833.3333
$ws.Range("I41").Value = 833.3333
$ws.Range("J41").Value = 0
$ws.Range("K41").Value = 2499.9999
$ws.Range("L41").Value = 0
$ws.Range("M41").Value = -2161.9999
$ws.Range("N41").Value = $null

$ws.Range("H131").Value = 2551.3333
$ws.Range("I131").Value = 3643
$ws.Range("J131").Value = 2005.5
$ws.Range("K131").Value = 10929
$ws.Range("L131").Value = 6016.5
$ws.Range("M131").Value = -5889
$ws.Range("N131").Value = -16096.5

$ws.Range("H133").Value = 1430
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").Value = $null

$ws.Range("H137").Value = 2414.5
$ws.Range("J137").Value = 3799
$ws.Range("L137").Value = 11397
$ws.Range("N137").Value = -21597

$ws.Range("H140").Value = 1835.3334
$ws.Range("I140").Value = 1835.3334
$ws.Range("K140").Value = 5506.0002
$ws.Range("M140").Value = -326.0002000000004

$ws.Range("H141").Value = 2998
$ws.Range("I141").Value = 2998
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 8994
$ws.Range("L141").Value = 0
$ws.Range("M141").Value = -3814
$ws.Range("N141").Value = $null

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 34000800
$ws.Range("I11").Value = 42500748
$ws.Range("J11").Value = 1004
$ws.Range("K11").Value = 42500748
$ws.Range("L11").Value = 1004
$ws.Range("M11").Value = -42500609
$ws.Range("N11").Value = -1282

$ws.Range("H14").Value = 1818750
$ws.Range("J14").Value = 25000
$ws.Range("L14").Value = 25000
$ws.Range("N14").Value = -25336

$ws.Range("H33").Value = 6000
$ws.Range("I33").Value = 0
$ws.Range("J33").Value = 6000
$ws.Range("K33").Value = 0
$ws.Range("L33").Value = 6000
$ws.Range("M33").Value = $null
$ws.Range("N33").Value = -6504

$ws.Range("H102").Value = 2295.4285
$ws.Range("I102").Value = 1811
$ws.Range("K102").Value = 1811
$ws.Range("M102").Value = -189

$ws.Range("H127").Value = 0
$ws.Range("J127").Value = 0
$ws.Range("L127").Value = 0
$ws.Range("N127").Value = $null

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H3").Value = 12800
$ws.Range("J3").Value = 12800
$ws.Range("L3").Value = 12800
$ws.Range("N3").Value = -13024

$ws.Range("H15").Value = 12800
$ws.Range("J15").Value = 12800
$ws.Range("L15").Value = 12800
$ws.Range("N15").Value = -13140

$ws.Range("H20").Value = 1050.909
$ws.Range("I20").Value = 1500
$ws.Range("J20").Value = 1006
$ws.Range("K20").Value = 1500
$ws.Range("L20").Value = 1006
$ws.Range("M20").Value = -1274
$ws.Range("N20").Value = -1458

$ws.Range("H22").Value = 966.6667
$ws.Range("I22").Value = 966.6667
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 966.6667
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -671.6667
$ws.Range("N22").Value = $null

$ws.Range("H27").Value = 966.6667
$ws.Range("I27").Value = 966.6667
$ws.Range("J27").Value = 0
$ws.Range("K27").Value = 966.6667
$ws.Range("L27").Value = 0
$ws.Range("M27").Value = -859.6667
$ws.Range("N27").Value = $null

$ws.Range("H46").Value = 0
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 0
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 0
$ws.Range("M46").Value = $null
$ws.Range("N46").Value = $null

$ws.Range("H55").Value = 1706.875
$ws.Range("I55").Value = 522.2857
$ws.Range("K55").Value = 522.2857
$ws.Range("M55").Value = -349.2857

$ws.Range("H100").Value = 2940.6
$ws.Range("I100").Value = 2940.6
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 2940.6
$ws.Range("L100").Value = 0
$ws.Range("M100").Value = -2399.6
$ws.Range("N100").Value = $null

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 7510002.5
$ws.Range("I14").Value = 15000000
$ws.Range("K14").Value = 15000000
$ws.Range("M14").Value = -14999832

$ws.Range("H51").Value = 2785
$ws.Range("I51").Value = 2785
$ws.Range("K51").Value = 2785
$ws.Range("M51").Value = -2275
